# Finished Week 13 logging
# Update the "H" row (row 2) target-depth counts on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 227
$wsOff.Range("C2").Value = 161
$wsOff.Range("D2").Value = 47
$wsOff.Range("E2").Value = 19
$wsOff.Range("F2").Value = 4

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 139
$wsDef.Range("C2").Value = 88
$wsDef.Range("D2").Value = 39
$wsDef.Range("E2").Value = 15
$wsDef.Range("F2").Value = 3
